# Update gh-pages output (generated at 456a3b4)
# Applies refreshed "想去人数" (interest count) numbers across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets, and inserts a newly announced
# performance (Ayasa LIVE TOUR) into 演出 ahead of the existing
# 花たん 2024 LIVE in Beijing row.

$wb = $excel.ActiveWorkbook

function Set-F {
    param($ws, [int]$row, $value)
    $ws.Cells.Item($row, 6).Value = $value
}

# ---- 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
Set-F $ws1 2  311
Set-F $ws1 4  457
Set-F $ws1 5  8767
Set-F $ws1 7  11177
Set-F $ws1 18 87
Set-F $ws1 22 1905
Set-F $ws1 23 711
Set-F $ws1 24 637
Set-F $ws1 25 360
Set-F $ws1 30 1305
Set-F $ws1 38 355
Set-F $ws1 39 309
Set-F $ws1 42 537
Set-F $ws1 43 391
Set-F $ws1 45 815
Set-F $ws1 46 658
Set-F $ws1 48 160
Set-F $ws1 49 147

# ---- 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
Set-F $ws2 4  23
Set-F $ws2 8  57
Set-F $ws2 18 67
Set-F $ws2 19 106

# Insert the new Ayasa LIVE TOUR row above row 24 (花たん shifts to row 25).
$ws2.Rows.Item(24).Insert()

function Set-Text {
    param($ws, [int]$row, [int]$col, [string]$value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-Text $ws2 24 2 "2024-11-24"
Set-Text $ws2 24 3 "北京·Ayasa LIVE TOUR 2024〜D.D.D.〜"
Set-Text $ws2 24 4 "建国门外郎家园10号61幢一层A3-06、二层A3-06号 EAST LIVE(东郎展演中心)"
Set-Text $ws2 24 5 "2024.11.24 19:00-11.24 20:30"
$ws2.Cells.Item(24, 6).Value = 49
$ws2.Cells.Item(24, 7).Value = 380
Set-Text $ws2 24 8 "https://show.bilibili.com/platform/detail.html?id=92778"
Set-Text $ws2 24 9 "//i0.hdslb.com/bfs/openplatform/202409/UQQJBlRE1727084578001.jpeg"

# Index column A keeps the plain running index (row-1); restore it/fix the
# now-shifted row below.
$ws2.Cells.Item(24, 1).Value = 23
$ws2.Cells.Item(25, 1).Value = 24

# ---- 本地生活 ----
$ws3 = $wb.Worksheets.Item("本地生活")
Set-F $ws3 3 2839
Set-F $ws3 4 350

# ---- 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
Set-F $ws4 2  311
Set-F $ws4 5  350
Set-F $ws4 7  8767
Set-F $ws4 9  11177
Set-F $ws4 19 1905
Set-F $ws4 20 711
Set-F $ws4 21 637
Set-F $ws4 22 360
Set-F $ws4 29 1305
Set-F $ws4 37 355
Set-F $ws4 40 537
Set-F $ws4 41 391
Set-F $ws4 46 658
Set-F $ws4 48 160
Set-F $ws4 49 147
